# [Fonds de solidarite] Add 2020-07-29 data
# Updates "nombre_aides" (column C) and "montant_total" (column D) figures
# for the rows affected by the 2020-07-29 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 9;  C = "46";  D = "113250.00" },
    @{ Row = 10; C = "309"; D = "938531.74" },
    @{ Row = 11; C = "127"; D = "500391.77" },
    @{ Row = 14; C = "13";  D = "34000.00" },
    @{ Row = 33; C = "448"; D = "1192724.87" },
    @{ Row = 34; C = "177"; D = "743340.11" },
    @{ Row = 35; C = "62";  D = "298974.00" },
    @{ Row = 36; C = "22";  D = "130500.00" },
    @{ Row = 49; C = "83";  D = "231937.17" },
    @{ Row = 50; C = "507"; D = "1575250.52" },
    @{ Row = 51; C = "221"; D = "853660.65" },
    @{ Row = 52; C = "72";  D = "380878.23" },
    @{ Row = 53; C = "22";  D = "121521.00" },
    @{ Row = 54; C = "15";  D = "42720.65" },
    @{ Row = 61; C = "206"; D = "463423.00" },
    @{ Row = 62; C = "26";  D = "57984.00" },
    @{ Row = 63; C = "183"; D = "450800.00" },
    @{ Row = 64; C = "96";  D = "274500.00" },
    @{ Row = 65; C = "32";  D = "123000.00" },
    @{ Row = 66; C = "14";  D = "28000.00" },
    @{ Row = 67; C = "83";  D = "208542.41" },
    @{ Row = 68; C = "352"; D = "1006792.54" },
    @{ Row = 69; C = "134"; D = "482662.18" },
    @{ Row = 70; C = "37";  D = "160849.00" },
    @{ Row = 71; C = "8";   D = "41000.00" },
    @{ Row = 72; C = "8";   D = "16000.00" },
    @{ Row = 73; C = "207"; D = "514326.09" },
    @{ Row = 74; C = "810"; D = "2371725.56" },
    @{ Row = 75; C = "297"; D = "1083266.79" },
    @{ Row = 76; C = "98";  D = "417484.52" },
    @{ Row = 78; C = "24";  D = "49500.00" }
)

foreach ($u in $updates) {
    # Force a text number format before assigning, otherwise Excel's COM
    # layer would auto-detect these numeric-looking strings and coerce
    # them into numeric cells (dropping trailing zeros, e.g. "113250.00"
    # -> 113250). The source column stores these as text values.
    $cCell = $ws.Range("C" + $u.Row)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C

    $dCell = $ws.Range("D" + $u.Row)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
}
